$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear cells that were removed entirely (C2, E2, C3)
$ws.Range("C2").Value = $null
$ws.Range("E2").Value = $null
$ws.Range("C3").Value = $null

# Update the slightly re-computed values (naive component forecaster bug fix)
$ws.Range("E3").Value = 2.807231216534301

$ws.Range("C4").Value = -0.9140166223623458
$ws.Range("E4").Value = 1.821983295885099

$ws.Range("C8").Value = -1.479696720105184

$ws.Range("E9").Value = -0.6155071485167807

$ws.Range("C11").Value = 2.192778679161966
$ws.Range("E11").Value = -0.5835597102573087

$ws.Range("C12").Value = 3.408364488606752

$ws.Range("E13").Value = 3.056075254340018

$ws.Range("C15").Value = 1.666553973046025
$ws.Range("E15").Value = -1.376301649685407

$ws.Range("C16").Value = 1.879266440112781
$ws.Range("E16").Value = -0.5015683214423916

$ws.Range("C17").Value = -2.620683231370935
$ws.Range("E17").Value = -3.531225750971467

$ws.Range("C18").Value = -3.036556262700263
